$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet tab
$ws.Name = "transitions"

# Merge the T_STR_ALT lexer state into T_STRING: update the one cell that
# still held the literal "T_STR_ALT" text so it now reads "T_STRING".
$ws.Range("K6").Value = "T_STRING"

# Update the view: unfreeze/refreeze pane at column B (was column R) and
# move the active selection to K6 (was AH15).
$ws.Range("B1").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("K6").Select()
